# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01293466051926884
$ws.Range("C2").Value = 0.04071648406533734
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("G2").Value = 1133.843310939226

$ws.Range("B3").Value = 0.6606524410359556
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 3.56341032713086
